# fix: export regvue json but some issue in parser with ipxact
#
# The "reserved" field appeared twice in the register maps (once for the
# high reserved bits, once for the low reserved bits) and both cells
# pointed at the same shared string. Split them into two distinct labels
# so downstream tooling (regvue / ipxact export) can tell them apart.

$wb = $excel.ActiveWorkbook

$block0 = $wb.Worksheets.Item("block0")
$block0.Range("C3").Value = "reserved1"
$block0.Range("C5").Value = "reserved0"

$block1 = $wb.Worksheets.Item("block1")
$block1.Range("C3").Value = "reserved1"
$block1.Range("C5").Value = "reserved0"

# Give block1's page a paper/orientation setup like the other sheets already have.
$block1.PageSetup.PaperSize = 9
$block1.PageSetup.Orientation = 1

# Move the current selection/active-tab: block0 becomes the active sheet
# (it was block1), with a fresh cursor position on each sheet.
$block0.Activate()
$block0.Range("I8").Select()

$block1.Activate()
$block1.Range("H10").Select()

$block0.Activate()
